# Revert "Added diary for last week":
# Clear out the 4 diary-entry rows (27-30) that were added for the
# previous week, restoring them to the same blank state as the rows
# that follow them (e.g. row 31), and update the view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Select()

# Copy the blank formatting from an already-empty row (31) onto rows
# 27:30 so the date column goes back to the "General" / unstyled look
# (style 13) instead of the date-formatted style (12) used for filled
# rows, then drop the values themselves.
$ws.Range("A31:G31").Copy() | Out-Null
$ws.Range("A27:G30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A27:G30").ClearContents()

# Restore the default (non-custom) row height now that the long
# wrapped diary text is gone.
$ws.Range("A27:G30").Rows.AutoFit() | Out-Null

# Restore the view/selection to where it was prior to the addition.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 23
$ws.Range("F26").Select() | Out-Null
